# Edit: PRUEBA_TEMPERATURA_6_22C.xlsx
# "Mas mediciones de temperatura" - replace raw temperature/timestamp dataset
# in "Datos crudos" with a newer run (34 samples instead of 37), update the
# lookup row pointer (H2), and let dependent formulas / chart caches refresh
# on recalculation + save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos crudos")

# New timestamp values (column C) for the 34 retained data rows.
$timestamps = @("2023-12-12 02:17:43","2023-12-12 02:18:44","2023-12-12 02:19:45","2023-12-12 02:20:46","2023-12-12 02:21:47","2023-12-12 02:22:48","2023-12-12 02:23:49","2023-12-12 02:24:50","2023-12-12 02:25:51","2023-12-12 02:26:52","2023-12-12 02:27:53","2023-12-12 02:28:54","2023-12-12 02:29:56","2023-12-12 02:30:57","2023-12-12 02:31:58","2023-12-12 02:32:59","2023-12-12 02:34:00","2023-12-12 02:35:01","2023-12-12 02:36:02","2023-12-12 02:37:03","2023-12-12 02:38:04","2023-12-12 02:39:05","2023-12-12 02:40:06","2023-12-12 02:41:08","2023-12-12 02:42:09","2023-12-12 02:43:10","2023-12-12 02:44:11","2023-12-12 02:45:12","2023-12-12 02:46:13","2023-12-12 02:47:14","2023-12-12 02:48:15","2023-12-12 02:49:16","2023-12-12 02:50:17","2023-12-12 02:51:18")

# New temperature readings (column E) for the same 34 rows.
$temperatures = @(26.590909090909001,26.4597902097901,25.541958041958001,25.4108391608391,24.886363636363601,24.493006993006901,24.493006993006901,23.968531468531399,24.0996503496503,24.0996503496503,23.837412587412501,23.837412587412501,23.575174825174798,23.444055944055901,23.7062937062937,23.444055944055901,23.575174825174798,23.312937062936999,23.312937062936999,23.181818181818102,23.181818181818102,23.181818181818102,23.181818181818102,23.0506993006993,23.312937062936999,23.181818181818102,23.181818181818102,22.919580419580399,23.444055944055901,23.181818181818102,23.181818181818102,22.919580419580399,23.181818181818102,23.0506993006993)

# Drop the last 3 raw samples (old rows 36-38) - the new run only has 34
# data rows instead of 37. Deleting shifts nothing else since there is no
# data below row 38.
$ws.Range("A36:A38").EntireRow.Delete()

# Overwrite the timestamp (C) and temperature (E) columns for rows 2..35
# with the new measurement run's data.
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $timestamps[$i]
    $ws.Cells.Item($row, 5).Value = $temperatures[$i]
}

# H2 is the manually-chosen "settling" row index fed into the
# INDEX(C:C,H2)-C2 settling-time formula in I2; bump it for the new run.
$ws.Range("H2").Value = 21
